$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 184, shifting existing rows 184:221 down to 185:222
$ws.Rows("184:184").Insert()

# Populate the newly inserted row 184 with the new record's data
$ws.Range("A184").Value = 3
$ws.Range("B184").Value = "Femacal de La Calera"
$ws.Range("C184").Value = "Coquimbo"
$ws.Range("D184").Value = 44476
$ws.Range("E184").Value = 5
$ws.Range("F184").Value = 100112032
$ws.Range("G184").Value = "Zapallo italiano"
$ws.Range("H184").Value = "Sin especificar"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 220
$ws.Range("K184").Value = 14000
$ws.Range("L184").Value = 15000
$ws.Range("M184").Value = 14500
$ws.Range("N184").Value = "`$/caja 70 unidades"
$ws.Range("O184").Value = "Región de Arica y Parinacota"
$ws.Range("P184").Value = 207
$ws.Range("Q184").Value = 70
$ws.Range("R184").Value = "Hortaliza"
